$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts B:F left to A:E
$ws.Range("A1").EntireColumn.Delete()
